$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.889.08'
$ws.Range('E2').Value = '  +1.09%  '
$ws.Range('D3').Value = '1.841.44'
$ws.Range('E3').Value = '  +1.12%  '
$ws.Range('D4').Value = '''1.008'
$ws.Range('E4').Value = '  -0.20%  '
$ws.Range('D5').Value = '''309.20'
$ws.Range('E5').Value = '  +1.28%  '
$ws.Range('D6').Value = '''1.006'
$ws.Range('E6').Value = '  -0.27%  '
$ws.Range('D7').Value = '''0.4747'
$ws.Range('E7').Value = '  +1.82%  '
$ws.Range('E8').Value = '  +2.31%  '
$ws.Range('D9').Value = '''0.07192'
$ws.Range('E9').Value = '  +1.04%  '
$ws.Range('D10').Value = '''0.9240'
$ws.Range('E10').Value = '  +2.77%  '
$ws.Range('E11').Value = '  +1.23%  '
$ws.Range('D12').Value = '''0.07643'
$ws.Range('E12').Value = '  -2.01%  '
$ws.Range('D13').Value = '1.902.20'
$ws.Range('E13').Value = '  +2.85%  '
$ws.Range('D14').Value = '''5.307'
$ws.Range('E14').Value = '  +1.12%  '
$ws.Range('D15').Value = '''6.399'
$ws.Range('E15').Value = '  +1.11%  '
$ws.Range('D16').Value = '''88.64'
$ws.Range('E17').Value = '  -0.23%  '
$ws.Range('D18').Value = '''0.000008632'
$ws.Range('E18').Value = '  +1.18%  '
$ws.Range('D19').Value = '''1.006'
$ws.Range('E19').Value = '  -0.28%  '
$ws.Range('D20').Value = '26.913.20'
$ws.Range('E20').Value = '  +1.00%  '
$ws.Range('D21').Value = '''14.53'
$ws.Range('E21').Value = '  +2.92%  '
$ws.Range('E22').Value = '  +0.80%  '
$ws.Range('E23').Value = '  +1.08%  '
$ws.Range('D24').Value = '''1.924'
$ws.Range('E24').Value = '  -0.60%  '
$ws.Range('D25').Value = '''152.21'
$ws.Range('E25').Value = '  +0.15%  '
$ws.Range('D26').Value = '''18.13'
$ws.Range('E26').Value = '  +1.43%  '
$ws.Range('D27').Value = '''2.001'
$ws.Range('E27').Value = '  +1.67%  '
$ws.Range('D28').Value = '''114.23'
$ws.Range('E28').Value = '  +0.61%  '
$ws.Range('D29').Value = '''4.940'
$ws.Range('E29').Value = '  +3.11%  '
$ws.Range('D30').Value = '''0.08849'
$ws.Range('E30').Value = '  +0.70%  '
$ws.Range('D31').Value = '''3.288'
$ws.Range('E31').Value = '  +5.27%  '
$ws.Range('D32').Value = '''0.7478'
$ws.Range('E32').Value = '  +2.74%  '
$ws.Range('D33').Value = '''1.170'
$ws.Range('E33').Value = '  +4.35%  '
$ws.Range('B34').Value = 'RenderToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D34').Value = '''2.753'
$ws.Range('E34').Value = '  +0.65%  '
$ws.Range('B35').Value = 'Filecoin'
$ws.Range('C35').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D35').Value = '''4.478'
$ws.Range('E35').Value = '  +1.09%  '
$ws.Range('E36').Value = '  +1.36%  '
$ws.Range('D37').Value = '''0.05257'
$ws.Range('E37').Value = '  +3.15%  '
$ws.Range('D38').Value = '''0.01946'
$ws.Range('E38').Value = '  +1.13%  '
$ws.Range('D39').Value = '''2.960'
$ws.Range('E39').Value = '  +1.36%  '
$ws.Range('D40').Value = '''0.5206'
$ws.Range('E40').Value = '  +3.70%  '
$ws.Range('D41').Value = '''6.959'
$ws.Range('E41').Value = '  +2.16%  '
$ws.Range('D42').Value = '''0.1512'
$ws.Range('E42').Value = '  +1.47%  '
$ws.Range('D43').Value = '''8.204'
$ws.Range('E43').Value = '  +3.08%  '
$ws.Range('D44').Value = '''10.53'
$ws.Range('E44').Value = '  +5.76%  '
$ws.Range('D45').Value = '''0.4718'
$ws.Range('E45').Value = '  +1.81%  '
$ws.Range('D46').Value = '''1.007'
$ws.Range('E46').Value = '  -0.28%  '
$ws.Range('D47').Value = '''101.53'
$ws.Range('E47').Value = '  +3.17%  '
$ws.Range('D48').Value = '''1.602'
$ws.Range('E48').Value = '  +3.36%  '
$ws.Range('D49').Value = '''65.43'
$ws.Range('E49').Value = '  +2.87%  '
$ws.Range('D50').Value = '''0.06026'
$ws.Range('E50').Value = '  +0.55%  '
$ws.Range('D51').Value = '''0.8839'
$ws.Range('E51').Value = '  +4.19%  '
